$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.117.78'
$ws.Range('E2').Value = '  +1.25%  '

$ws.Range('D3').Value = '3.503.73'
$ws.Range('E3').Value = '  +0.28%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '603.16'
$ws.Range('E5').Value = '  +0.48%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.32'
$ws.Range('E6').Value = '  +3.40%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.609'
$ws.Range('E7').Value = '  -0.91%  '

$ws.Range('D8').Value = '3.496.11'
$ws.Range('E8').Value = '  +0.27%  '

$ws.Range('E10').Value = '  -0.13%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.23'
$ws.Range('E11').Value = '  +9.22%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.580'
$ws.Range('E12').Value = '  +0.59%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '46.05'
$ws.Range('E13').Value = '  -1.57%  '

$ws.Range('E14').Value = '  -0.60%  '

$ws.Range('D15').Value = '4.070.11'
$ws.Range('E15').Value = '  +0.33%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.26'
$ws.Range('E16').Value = '  -0.08%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '608.93'
$ws.Range('E17').Value = '  +0.06%  '

$ws.Range('D18').Value = '3.509.05'
$ws.Range('E18').Value = '  +0.29%  '

$ws.Range('D19').Value = '70.253.80'
$ws.Range('E19').Value = '  +1.36%  '

$ws.Range('E20').Value = '  +0.84%  '

$ws.Range('E21').Value = '  +0.95%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.872'
$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('E23').Value = '  -14.83%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '15.50'
$ws.Range('E24').Value = '  -0.90%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '97.13'
$ws.Range('E25').Value = '  +1.58%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.71'
$ws.Range('E26').Value = '  -2.98%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('E28').Value = '  -1.69%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.65'
$ws.Range('E29').Value = '  +2.35%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.93'
$ws.Range('E30').Value = '  -2.92%  '

$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.96'
$ws.Range('E31').Value = '  -3.83%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.00'
$ws.Range('E32').Value = '  -4.60%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '643.81'
$ws.Range('E33').Value = '  +15.97%  '

$ws.Range('E34').Value = '  -3.93%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.81'
$ws.Range('E35').Value = '  +0.00%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.54'
$ws.Range('E36').Value = '  +2.95%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0988'
$ws.Range('E37').Value = '  -1.80%  '

$ws.Range('E38').Value = '  +0.02%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0470'
$ws.Range('E39').Value = '  +5.13%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '56.56'
$ws.Range('E40').Value = '  -0.20%  '

$ws.Range('E41').Value = '  +0.07%  '

$ws.Range('E42').Value = '  +1.52%  '

$ws.Range('D43').Value = '3.349.69'
$ws.Range('E43').Value = '  -0.31%  '

$ws.Range('D44').Value = '0.0₃0733'
$ws.Range('E44').Value = '  +6.12%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.306'
$ws.Range('E45').Value = '  -5.10%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '32.11'
$ws.Range('E46').Value = '  -1.82%  '

$ws.Range('E47').Value = '  +1.20%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.53'
$ws.Range('E48').Value = '  -2.16%  '

$ws.Range('E49').Value = '  +0.61%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '133.76'
$ws.Range('E50').Value = '  -0.28%  '

$ws.Range('E51').Value = '  -0.01%  '
